$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for 4 new transaction rows just above the (blank gap row
# that precedes the) GRAND TOTAL row - this shifts GRAND TOTAL and the
# footer rows down by 4, carrying their formatting/merge with them.
$ws.Range("A11:H14").Insert()

# Row 11: transaction #8
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "Baterai ABC AA"
$ws.Range("D11").Value = "Elektronik"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 6000
$ws.Range("G11").Value = 6000
$ws.Range("H11").Value = "01/07/2025 00:06:44"

# Row 12: transaction #9
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Air Mineral Aqua 600ml"
$ws.Range("D12").Value = "Minuman"
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 3500
$ws.Range("G12").Value = 17500
$ws.Range("H12").Value = "01/07/2025 00:26:44"

# Row 13: transaction #10
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = "Ayam Geprek"
$ws.Range("D13").Value = "Depi"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 20000
$ws.Range("G13").Value = 20000
$ws.Range("H13").Value = "01/07/2025 10:09:05"

# Row 14: transaction #11
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 11
$ws.Range("C14").Value = "a"
$ws.Range("D14").Value = "a"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1111
$ws.Range("G14").Value = 1111
$ws.Range("H14").Value = "02/07/2025 23:28:07"

# Row 16: GRAND TOTAL now covers all 11 transactions
$ws.Range("G16").Value = 223111

# Footer rows (18, 19, 20)
$ws.Range("A18").Value = "Total Transaksi: 11"
$ws.Range("A19").Value = "Periode: Mingguan - 30/06/2025 s/d 06/07/2025"
$ws.Range("A20").Value = "Dibuat pada: 02/07/2025 23:59:17"
